$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "Importe" cells (H2:H4) are stored as text (shared strings) that
# use a comma-based (es-AR) decimal/thousands format, e.g. "98,00",
# "9.174,00", "48.000,00". Fix the scraped formatting so the numbers use a
# plain dot-decimal, no-thousands-separator textual representation instead
# (98.00 / 9174.00 / 48000.00), while keeping them as text values.

# Force a Text number format first so Excel does not reinterpret the
# replaced values as numeric when they look like plain numbers.
$ws.Range("H2:H4").NumberFormat = "@"

$ws.Cells.Replace("98,00", "98.00", 1)
$ws.Cells.Replace("9.174,00", "9174.00", 1)
$ws.Cells.Replace("48.000,00", "48000.00", 1)
